$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '292.15'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.35%'

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '31.05'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '0.76%'

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.975'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.64%'

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07460'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '2.96%'

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.240'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-5.04%'

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '7.761'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.11%'

# Row 8
$ws.Range("B8").Value = 'GateToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.776'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '1.85%'

# Row 9
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9188'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '2.50%'

# Row 10
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09357'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '18.15%'

# Row 11
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1734'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '4.09%'

# Row 12
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08341'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '2.76%'

# Row 13
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03277'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '5.47%'

# Row 14
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09950'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.80%'

# Row 15
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001503'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.41%'

# Row 16
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005717'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-2.49%'

# Row 17
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.467'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.06%'

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.155'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '3.88%'

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3329'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.31%'

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1308'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.77%'

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.103'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '3.42%'

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2125'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-3.48%'

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04532'

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001221'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.91%'

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004260'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-3.46%'

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001300'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.03%'

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0003396'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '0.03%'

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01618'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '3.07%'

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04569'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '4.34%'

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007453'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '1.96%'

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009851'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '0.18%'

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1357'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '3.42%'

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002160'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '7.53%'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.008953'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-5.68%'

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006105'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '6.00%'

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000750'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.05%'

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.794'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '24.70%'

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002001'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-30.99%'

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002100'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.05%'

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002000'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.05%'
